$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update non-price columns (B, C, E) with new text values
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +7.43%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("E12").Value = "  +6.97%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("E15").Value = "  +8.50%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +32.18%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  +12.06%  "
$ws.Range("E26").Value = "  +6.24%  "
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("E30").Value = "  +7.05%  "
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E40").Value = "  +10.86%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("E41").Value = "  +24.31%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +4.87%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  -0.13%  "

# Update price column (D) values as text, forcing text format to avoid
# Excel auto-converting numeric-looking strings into numbers
$dCells = @("D2","D3","D5","D6","D7","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D24","D27","D28","D30","D31","D32","D33","D34","D37","D38","D40","D41","D42","D45","D47","D48","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D2").Value = "37.119.16"
$ws.Range("D3").Value = "2.054.27"
$ws.Range("D5").Value = "249.82"
$ws.Range("D6").Value = "0.672"
$ws.Range("D7").Value = "59.64"
$ws.Range("D9").Value = "0.391"
$ws.Range("D10").Value = "0.0794"
$ws.Range("D12").Value = "16.08"
$ws.Range("D13").Value = "2.355.65"
$ws.Range("D14").Value = "0.841"
$ws.Range("D15").Value = "5.76"
$ws.Range("D16").Value = "2.056.27"
$ws.Range("D17").Value = "18.91"
$ws.Range("D18").Value = "37.091.40"
$ws.Range("D19").Value = "75.86"
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("D21").Value = "5.46"
$ws.Range("D22").Value = "238.86"
$ws.Range("D24").Value = "2.43"
$ws.Range("D27").Value = "169.62"
$ws.Range("D28").Value = "20.33"
$ws.Range("D30").Value = "1.13"
$ws.Range("D31").Value = "4.82"
$ws.Range("D32").Value = "0.0633"
$ws.Range("D33").Value = "4.55"
$ws.Range("D34").Value = "0.0898"
$ws.Range("D37").Value = "1.74"
$ws.Range("D38").Value = "0.107"
$ws.Range("D40").Value = "3.11"
$ws.Range("D41").Value = "5.16"
$ws.Range("D42").Value = "17.77"
$ws.Range("D45").Value = "97.96"
$ws.Range("D47").Value = "1.297.09"
$ws.Range("D48").Value = "3.85"
$ws.Range("D49").Value = "2.88"
$ws.Range("D51").Value = "2.240.01"
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
